# Mise à jour de l'application
# Add two new attendance-tracking columns (DL, DM) for 2026-02-05 and
# 2026-02-06, mirroring the formatting of the last existing date column
# (DK) and filling in each player's status for the two new sessions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone the formatting of the last date column (DK, rows 1-31) onto the
#    two new columns (DL, DM) so number formats / styles match exactly.
$ws.Range("DK1:DK31").Copy()
$ws.Range("DL1:DM31").PasteSpecial(-4122)

# Rows 12, 21 and 23 stop well before column DK in the source sheet (their
# data doesn't extend that far), so they must NOT get new DL/DM cells.
$ws.Range("DL12:DM12").Clear()
$ws.Range("DL21:DM21").Clear()
$ws.Range("DL23:DM23").Clear()

# 2) New date headers for row 1.
$ws.Range("DL1").Value = 46058
$ws.Range("DM1").Value = 46059

# 3) Attendance marks for the two new sessions, per player row.
$ws.Range("DL2").Value = "P"
$ws.Range("DM2").Value = "P"

$ws.Range("DL3").Value = "P"
$ws.Range("DM3").Value = "R"

$ws.Range("DL4").Value = "P"
$ws.Range("DM4").Value = "P"

$ws.Range("DL5").Value = "P"
$ws.Range("DM5").Value = "P"

$ws.Range("DL6").Value = "B"
$ws.Range("DM6").Value = "B"

$ws.Range("DL7").Value = "P"
$ws.Range("DM7").Value = "P"

$ws.Range("DL8").Value = "P"
$ws.Range("DM8").Value = "P"

$ws.Range("DL9").Value = "P"
$ws.Range("DM9").Value = "P"

$ws.Range("DL10").Value = "P"
$ws.Range("DM10").Value = "P"

$ws.Range("DL11").Value = "P"
$ws.Range("DM11").Value = "P"

# Row 12 (Yanis Berrached) is left untouched - no DL/DM cells.

$ws.Range("DL13").Value = "P"
$ws.Range("DM13").Value = "REP"

$ws.Range("DL14").Value = "P"
$ws.Range("DM14").Value = "P"

$ws.Range("DL15").Value = "P"
$ws.Range("DM15").Value = "P"

# Rows 16 & 17 get formatted-but-empty DL/DM cells (already handled above).

$ws.Range("DL18").Value = "B"
$ws.Range("DM18").Value = "B"

$ws.Range("DL19").Value = "P"
$ws.Range("DM19").Value = "P"

$ws.Range("DL20").Value = "P"
$ws.Range("DM20").Value = "P"

# Row 21 (Amir Kherrab) is left untouched - no DL/DM cells.

$ws.Range("DL22").Value = "P"
$ws.Range("DM22").Value = "P"

# Row 23 (Wael Fareh) is left untouched - no DL/DM cells.

$ws.Range("DL24").Value = "P"
$ws.Range("DM24").Value = "P"

# Row 25 gets formatted-but-empty DL/DM cells (already handled above).

$ws.Range("DL26").Value = "P"
$ws.Range("DM26").Value = "P"

$ws.Range("DL27").Value = "P"
$ws.Range("DM27").Value = "P"

$ws.Range("DL28").Value = "P"
$ws.Range("DM28").Value = "P"

$ws.Range("DL29").Value = "P"
$ws.Range("DM29").Value = "P"

$ws.Range("DL30").Value = "P"
$ws.Range("DM30").Value = "P"

$ws.Range("DL31").Value = "P"
$ws.Range("DM31").Value = "P"

# 4) Update the view state to match: freeze pane scrolled to show the new
#    columns, with the active selection on DO18 (matches the post-edit
#    selection in the saved file).
$ws.Range("DO18").Select()

# 5) Recalculate so the COUNTA/COUNTIF summary formulas (columns B-J) pick
#    up the two new columns of data.
$excel.Calculate()
